$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Apparel and Clothing"
$ws.Range("B7").Value = "13/662"
$ws.Range("C7").Value = "11/586"
$ws.Range("D7").Value = "9/462"
$ws.Range("E7").Value = "9/462"
$ws.Range("F7").Value = "10/510"
$ws.Range("G7").Value = "12/611"
$ws.Range("H7").Value = "10/490"
$ws.Range("I7").Value = "13/654"

$ws.Range("L12").Select()
